$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Correct the DAQ codes in column D (rows 2-25) ---
$daqCodes = @{
    2  = "000001"
    3  = "100001"
    4  = "010001"
    5  = "110001"
    6  = "001001"
    7  = "101001"
    8  = "011001"
    9  = "111001"
    10 = "000010"
    11 = "100010"
    12 = "010010"
    13 = "110010"
    14 = "001010"
    15 = "101010"
    16 = "011010"
    17 = "111010"
    18 = "000100"
    19 = "100100"
    20 = "010100"
    21 = "110100"
    22 = "001100"
    23 = "101100"
    24 = "011100"
    25 = "111100"
}

foreach ($row in $daqCodes.Keys) {
    $ws.Cells.Item($row, 4).Value = $daqCodes[$row]
}

# --- New working column F (rows 2-25), formatted like column D ---
$fRange = $ws.Range("F2:F25")
$fRange.NumberFormat = "@"
$fRange.HorizontalAlignment = -4108

# --- Update the active selection left by the author ---
$ws.Range("H7").Select()
